# Applies the CATALOG_sample.xlsx edit described by the commit:
# "Adding examples, changing URIs to https"
#
# Summary of changes:
#  - Collection sheet: RELATION:Contact -> RELATION:contactPoint, value
#    changed from a plain name to an email address, a mailto: hyperlink is
#    added on B6, and the old review comment on B6 is removed (which also
#    drops the now-unused legacyDrawing/vmlDrawing + comments part).
#  - People sheet: a new example "contactPoint" column (I) is added.
#  - Actions sheet: a new "endTime" column (D) is inserted (shifting the
#    following RELATION columns right), with example end-time and
#    ContentLocation data filled in.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Collection sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Collection")

# Remove the legacy reviewer comment that lived on B6 ("Put the ID of a
# Person or Organisation here"). Deleting it also drops the vmlDrawing /
# legacyDrawing relationship that only existed to support it.
$comment = $ws1.Range("B6").Comment
if ($comment -ne $null) {
    $comment.Delete()
}

$ws1.Range("A6").Value = "RELATION:contactPoint"
$ws1.Range("B6").Value = "peter.sefton@uts.edu.au"

# Link the new email value.
$ws1.Hyperlinks.Add($ws1.Range("B6"), "mailto:peter.sefton@uts.edu.au")
# Adding a hyperlink auto-applies Excel's built-in "Hyperlink" style;
# the source value was plain text, so restore the default cell style.
$ws1.Range("B6").Style = "Normal"

# ---------------------------------------------------------------------
# People sheet - add an example contactPoint column
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("People")

$ws3.Range("I1").Value = "contactPoint>TYPE:ContactPoint>"
$ws3.Range("I2").Value = "ID: peter.sefton@uts.edu.au, contactType: customer service, email: peter.sefton@uts.edu.au, URL: http://orcid.org/0000-0002-3545-944X, name: Contact Peter Sefton"

# ---------------------------------------------------------------------
# Actions sheet - insert a new "endTime" column before RELATION:Result
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Actions")

$ws6.Columns.Item(4).Insert()

$ws6.Range("D1").Value = "endTime"
$ws6.Range("D2").Value = "2017:06:11T12:56:14+10:00"
$ws6.Range("F2").Value = "Catalina Park"
$ws6.Range("D3").Value = "2018:09:19T17:01:07+10:00"
